$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf18"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.350653
$ws.Range("H2").Value = 1.051959
$ws.Range("I2").Value = 0.04536179359243143
$ws.Range("J2").Value = 0.04536179359243143
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.714474
$ws.Range("N2").Value = 2.143422
$ws.Range("O2").Value = 0.138796410342318
$ws.Range("P2").Value = 0.138796410342318
$ws.Range("Q2").Value = 0.250532451522
$ws.Range("R2").Value = 2.254792063698
$ws.Range("S2").Value = 0.006296054117318646
$ws.Range("T2").Value = 0.006296054117318646

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf18"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.350653
$ws.Range("H3").Value = 1.051959
$ws.Range("I3").Value = 0.04536179359243143
$ws.Range("J3").Value = 0.04536179359243143
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 4.140873
$ws.Range("N3").Value = 12.422619
$ws.Range("O3").Value = 0.8044215857867821
$ws.Range("P3").Value = 0.8044215857867821
$ws.Range("Q3").Value = 1.452009540069
$ws.Range("R3").Value = 13.068085860621
$ws.Range("S3").Value = 0.03649000593575639
$ws.Range("T3").Value = 0.03649000593575639

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf18"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.350653
$ws.Range("H4").Value = 1.051959
$ws.Range("I4").Value = 0.04536179359243143
$ws.Range("J4").Value = 0.04536179359243143
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.2847646666666667
$ws.Range("N4").Value = 0.8542940000000001
$ws.Range("O4").Value = 0.05531945672713084
$ws.Range("P4").Value = 0.05531945672713083
$ws.Range("Q4").Value = 0.0998535846606667
$ws.Range("R4").Value = 0.8986822619460002
$ws.Range("S4").Value = 0.002509389777701552
$ws.Range("T4").Value = 0.002509389777701551

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf18"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.350653
$ws.Range("H5").Value = 1.051959
$ws.Range("I5").Value = 0.04536179359243143
$ws.Range("J5").Value = 0.04536179359243143
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.007528666666666667
$ws.Range("N5").Value = 0.022586
$ws.Range("O5").Value = 0.00146254714376898
$ws.Range("P5").Value = 0.00146254714376898
$ws.Range("Q5").Value = 0.002639949552666667
$ws.Range("R5").Value = 0.02375954597400001
$ws.Range("S5").Value = 0.00006634376165484862
$ws.Range("T5").Value = 0.00006634376165484861

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf18"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 6.557365999999999
$ws.Range("H6").Value = 19.672098
$ws.Range("I6").Value = 0.8482855786262421
$ws.Range("J6").Value = 0.8482855786262421
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.714474
$ws.Range("N6").Value = 2.143422
$ws.Range("O6").Value = 0.138796410342318
$ws.Range("P6").Value = 0.138796410342318
$ws.Range("Q6").Value = 4.685067515484
$ws.Range("R6").Value = 42.165607639356
$ws.Range("S6").Value = 0.1177389932584786
$ws.Range("T6").Value = 0.1177389932584786

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf18"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 6.557365999999999
$ws.Range("H7").Value = 19.672098
$ws.Range("I7").Value = 0.8482855786262421
$ws.Range("J7").Value = 0.8482855786262421
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 4.140873
$ws.Range("N7").Value = 12.422619
$ws.Range("O7").Value = 0.8044215857867821
$ws.Range("P7").Value = 0.8044215857867821
$ws.Range("Q7").Value = 27.153219820518
$ws.Range("R7").Value = 244.378978384662
$ws.Range("S7").Value = 0.6823792303585797
$ws.Range("T7").Value = 0.6823792303585797

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf18"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 6.557365999999999
$ws.Range("H8").Value = 19.672098
$ws.Range("I8").Value = 0.8482855786262421
$ws.Range("J8").Value = 0.8482855786262421
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 0.2847646666666667
$ws.Range("N8").Value = 0.8542940000000001
$ws.Range("O8").Value = 0.05531945672713084
$ws.Range("P8").Value = 0.05531945672713083
$ws.Range("Q8").Value = 1.867306143201333
$ws.Range("R8").Value = 16.805755288812
$ws.Range("S8").Value = 0.04692669735906355
$ws.Range("T8").Value = 0.04692669735906354

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf18"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 6.557365999999999
$ws.Range("H9").Value = 19.672098
$ws.Range("I9").Value = 0.8482855786262421
$ws.Range("J9").Value = 0.8482855786262421
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.007528666666666667
$ws.Range("N9").Value = 0.022586
$ws.Range("O9").Value = 0.00146254714376898
$ws.Range("P9").Value = 0.00146254714376898
$ws.Range("Q9").Value = 0.04936822282533333
$ws.Range("R9").Value = 0.444314005428
$ws.Range("S9").Value = 0.001240657650120227
$ws.Range("T9").Value = 0.001240657650120227

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf18"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.7457606666666666
$ws.Range("H10").Value = 2.237282
$ws.Range("I10").Value = 0.09647441040198541
$ws.Range("J10").Value = 0.09647441040198541
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.714474
$ws.Range("N10").Value = 2.143422
$ws.Range("O10").Value = 0.138796410342318
$ws.Range("P10").Value = 0.138796410342318
$ws.Range("Q10").Value = 0.532826606556
$ws.Range("R10").Value = 4.795439459004
$ws.Range("S10").Value = 0.01339030185368716
$ws.Range("T10").Value = 0.01339030185368716

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Fgf18"
$ws.Range("C11").Value = "Fgfr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 0.7457606666666666
$ws.Range("H11").Value = 2.237282
$ws.Range("I11").Value = 0.09647441040198541
$ws.Range("J11").Value = 0.09647441040198541
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 4.140873
$ws.Range("N11").Value = 12.422619
$ws.Range("O11").Value = 0.8044215857867821
$ws.Range("P11").Value = 0.8044215857867821
$ws.Range("Q11").Value = 3.088100209062
$ws.Range("R11").Value = 27.792901881558
$ws.Range("S11").Value = 0.07760609820340994
$ws.Range("T11").Value = 0.07760609820340994

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Fgf18"
$ws.Range("C12").Value = "Fgfr2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.7457606666666666
$ws.Range("H12").Value = 2.237282
$ws.Range("I12").Value = 0.09647441040198541
$ws.Range("J12").Value = 0.09647441040198541
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.2847646666666667
$ws.Range("N12").Value = 0.8542940000000001
$ws.Range("O12").Value = 0.05531945672713084
$ws.Range("P12").Value = 0.05531945672713083
$ws.Range("Q12").Value = 0.2123662876564445
$ws.Range("R12").Value = 1.911296588908
$ws.Range("S12").Value = 0.005336911971508094
$ws.Range("T12").Value = 0.005336911971508092

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Fgf18"
$ws.Range("C13").Value = "Fgfr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.7457606666666666
$ws.Range("H13").Value = 2.237282
$ws.Range("I13").Value = 0.09647441040198541
$ws.Range("J13").Value = 0.09647441040198541
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.007528666666666667
$ws.Range("N13").Value = 0.022586
$ws.Range("O13").Value = 0.00146254714376898
$ws.Range("P13").Value = 0.00146254714376898
$ws.Range("Q13").Value = 0.005614583472444444
$ws.Range("R13").Value = 0.050531251252
$ws.Range("S13").Value = 0.0001410983733802201
$ws.Range("T13").Value = 0.0001410983733802201

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Fgf18"
$ws.Range("C14").Value = "Fgfr2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1.0
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.07636
$ws.Range("H14").Value = 0.22908
$ws.Range("I14").Value = 0.009878217379341012
$ws.Range("J14").Value = 0.009878217379341012
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 0.714474
$ws.Range("N14").Value = 2.143422
$ws.Range("O14").Value = 0.138796410342318
$ws.Range("P14").Value = 0.138796410342318
$ws.Range("Q14").Value = 0.05455723464
$ws.Range("R14").Value = 0.4910151117600001
$ws.Range("S14").Value = 0.001371061112833632
$ws.Range("T14").Value = 0.001371061112833632

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Fgf18"
$ws.Range("C15").Value = "Fgfr2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1.0
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.07636
$ws.Range("H15").Value = 0.22908
$ws.Range("I15").Value = 0.009878217379341012
$ws.Range("J15").Value = 0.009878217379341012
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 4.140873
$ws.Range("N15").Value = 12.422619
$ws.Range("O15").Value = 0.8044215857867821
$ws.Range("P15").Value = 0.8044215857867821
$ws.Range("Q15").Value = 0.31619706228
$ws.Range("R15").Value = 2.84577356052
$ws.Range("S15").Value = 0.007946251289036048
$ws.Range("T15").Value = 0.007946251289036048

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Fgf18"
$ws.Range("C16").Value = "Fgfr2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1.0
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.07636
$ws.Range("H16").Value = 0.22908
$ws.Range("I16").Value = 0.009878217379341012
$ws.Range("J16").Value = 0.009878217379341012
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 0.2847646666666667
$ws.Range("N16").Value = 0.8542940000000001
$ws.Range("O16").Value = 0.05531945672713084
$ws.Range("P16").Value = 0.05531945672713083
$ws.Range("Q16").Value = 0.02174462994666667
$ws.Range("R16").Value = 0.19570166952
$ws.Range("S16").Value = 0.000546457618857647
$ws.Range("T16").Value = 0.0005464576188576469

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Fgf18"
$ws.Range("C17").Value = "Fgfr2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1.0
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.07636
$ws.Range("H17").Value = 0.22908
$ws.Range("I17").Value = 0.009878217379341012
$ws.Range("J17").Value = 0.009878217379341012
$ws.Range("K17").Value = 2.0
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.007528666666666667
$ws.Range("N17").Value = 0.022586
$ws.Range("O17").Value = 0.00146254714376898
$ws.Range("P17").Value = 0.00146254714376898
$ws.Range("Q17").Value = 0.0005748889866666666
$ws.Range("R17").Value = 0.005174000880000001
$ws.Range("S17").Value = 0.0000144473586136843
$ws.Range("T17").Value = 0.00001444735861368429
